$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.897.58"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.48"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.91"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3794"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08899"
$ws.Range("E9").Value = "  -9.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.50"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.286"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.862.82"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.195"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06628"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.83"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.078"
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.951.78"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.285"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.087.78"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.503"
$ws.Range("E27").Value = "  -6.78%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.90"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.62"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.78"
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1051"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.554"
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.592"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.272"
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06503"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02390"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2173"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.264"
$ws.Range("E39").Value = "  +6.52%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.193"
$ws.Range("E40").Value = "  -6.19%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.69"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6336"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.872"
$ws.Range("E43").Value = "  -4.41%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5961"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.07"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.280"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.665"
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.959"
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.50"
$ws.Range("E51").Value = "  -3.61%  "
